# Scheduled runner update: refresh computed price/profit figures on the
# per-job Leve sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR) to reflect the
# latest market data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 20216.666
$ws.Range("J82").Value = 29500
$ws.Range("L82").Value = 88500
$ws.Range("N82").Value = -89312

$ws.Range("H85").Value = 20216.666
$ws.Range("J85").Value = 29500
$ws.Range("L85").Value = 88500
$ws.Range("N85").Value = -91308

$ws.Range("H96").Value = 654.6667

$ws.Range("H101").Value = 2262.6
$ws.Range("I101").Value = 2262.6
$ws.Range("K101").Value = 6787.799999999999
$ws.Range("M101").Value = -5165.799999999999

$ws.Range("H107").Value = 740.6667
$ws.Range("I107").Value = 885.6799999999999
$ws.Range("K107").Value = 885.6799999999999
$ws.Range("M107").Value = 1034.32

$ws.Range("H113").Value = 3565.5557
$ws.Range("I113").Value = 2499.4443
$ws.Range("J113").Value = 4631.6665
$ws.Range("K113").Value = 2499.4443
$ws.Range("L113").Value = 4631.6665
$ws.Range("M113").Value = 754.5556999999999
$ws.Range("N113").Value = -11139.6665

$ws.Range("H132").Value = 4475
$ws.Range("I132").Value = 4173.3335
$ws.Range("J132").Value = 4977.778
$ws.Range("K132").Value = 12520.0005
$ws.Range("L132").Value = 14933.334
$ws.Range("M132").Value = -9990.000499999998
$ws.Range("N132").Value = -19993.334

$ws.Range("H137").Value = 2067.5454
$ws.Range("I137").Value = 1510.409
$ws.Range("J137").Value = 3181.818
$ws.Range("K137").Value = 4531.227000000001
$ws.Range("L137").Value = 9545.454000000002
$ws.Range("M137").Value = -1981.227000000001
$ws.Range("N137").Value = -14645.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 14474.25
$ws.Range("J55").Value = 14474.25
$ws.Range("L55").Value = 14474.25
$ws.Range("N55").Value = -15104.25

$ws.Range("H80").Value = 19757.334
$ws.Range("J80").Value = 19757.334
$ws.Range("L80").Value = 19757.334
$ws.Range("N80").Value = -21753.334

$ws.Range("H83").Value = 19757.334
$ws.Range("J83").Value = 19757.334
$ws.Range("L83").Value = 59272.00199999999
$ws.Range("N83").Value = -69256.00199999999

$ws.Range("H122").Value = 2006.1111
$ws.Range("I122").Value = 1916
$ws.Range("K122").Value = 5748
$ws.Range("M122").Value = -3298

$ws.Range("H132").Value = 1965.5714
$ws.Range("I132").Value = 1710.4333
$ws.Range("J132").Value = 3496.4
$ws.Range("K132").Value = 5131.2999
$ws.Range("L132").Value = 10489.2
$ws.Range("M132").Value = -2601.2999
$ws.Range("N132").Value = -15549.2

$ws.Range("H139").Value = 108928.57
$ws.Range("J139").Value = 122083.336
$ws.Range("L139").Value = 122083.336
$ws.Range("N139").Value = -132363.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 35000
$ws.Range("J35").Value = 35000
$ws.Range("L35").Value = 35000
$ws.Range("N35").Value = -35620

$ws.Range("H51").Value = 34000
$ws.Range("J51").Value = 34000
$ws.Range("L51").Value = 34000
$ws.Range("N51").Value = -34982

$ws.Range("H82").Value = 23966
$ws.Range("J82").Value = 29834.47
$ws.Range("L82").Value = 29834.47
$ws.Range("N82").Value = -30600.47

$ws.Range("H85").Value = 23966
$ws.Range("J85").Value = 29834.47
$ws.Range("L85").Value = 29834.47
$ws.Range("N85").Value = -32486.47

$ws.Range("H94").Value = 443.4762
$ws.Range("I94").Value = 367.375
$ws.Range("J94").Value = 687
$ws.Range("K94").Value = 367.375
$ws.Range("L94").Value = 687
$ws.Range("M94").Value = 83.625
$ws.Range("N94").Value = -1589

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6665
$ws.Range("I16").Value = 7600
$ws.Range("J16").Value = 1990
$ws.Range("K16").Value = 7600
$ws.Range("L16").Value = 1990
$ws.Range("M16").Value = -7313
$ws.Range("N16").Value = -2564

$ws.Range("H31").Value = 2744.6562
$ws.Range("I31").Value = 2027.5834
$ws.Range("J31").Value = 3174.9
$ws.Range("K31").Value = 2027.5834
$ws.Range("L31").Value = 3174.9
$ws.Range("M31").Value = -1732.5834
$ws.Range("N31").Value = -3764.9

$ws.Range("H34").Value = 2744.6562
$ws.Range("I34").Value = 2027.5834
$ws.Range("J34").Value = 3174.9
$ws.Range("K34").Value = 2027.5834
$ws.Range("L34").Value = 3174.9
$ws.Range("M34").Value = -1825.5834
$ws.Range("N34").Value = -3578.9

$ws.Range("H41").Value = 15907
$ws.Range("J41").Value = 19758.75
$ws.Range("L41").Value = 19758.75
$ws.Range("N41").Value = -20614.75

$ws.Range("H50").Value = 8915.143
$ws.Range("J50").Value = 8915.143
$ws.Range("L50").Value = 8915.143
$ws.Range("N50").Value = -10165.143

$ws.Range("H107").Value = 1054.3846
$ws.Range("I107").Value = 1210.2106
$ws.Range("J107").Value = 631.4286
$ws.Range("K107").Value = 1210.2106
$ws.Range("L107").Value = 631.4286
$ws.Range("M107").Value = 709.7893999999999
$ws.Range("N107").Value = -4471.4286

$ws.Range("H109").Value = 10966.667
$ws.Range("J109").Value = 10966.667
$ws.Range("L109").Value = 10966.667
$ws.Range("N109").Value = -13046.667

$ws.Range("H113").Value = 6665
$ws.Range("I113").Value = 7600
$ws.Range("J113").Value = 1990
$ws.Range("K113").Value = 7600
$ws.Range("L113").Value = 1990
$ws.Range("M113").Value = -5430
$ws.Range("N113").Value = -6330

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 34500
$ws.Range("J57").Value = 34500
$ws.Range("L57").Value = 34500
$ws.Range("N57").Value = -36140

$ws.Range("H102").Value = 3067.7334
$ws.Range("I102").Value = 1922
$ws.Range("K102").Value = 1922
$ws.Range("M102").Value = -300

$ws.Range("H113").Value = 1594.88
$ws.Range("I113").Value = 1536.0667
$ws.Range("J113").Value = 1683.1
$ws.Range("K113").Value = 1536.0667
$ws.Range("L113").Value = 1683.1
$ws.Range("M113").Value = 633.9332999999999
$ws.Range("N113").Value = -6023.1

$ws.Range("H122").Value = 2334.111
$ws.Range("I122").Value = 2250.875
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6752.625
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4302.625
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1810
$ws.Range("I40").Value = 1637.5
$ws.Range("K40").Value = 1637.5
$ws.Range("M40").Value = -1501.5

$ws.Range("H61").Value = 2513.8667
$ws.Range("I61").Value = 1770.8
$ws.Range("K61").Value = 1770.8
$ws.Range("M61").Value = -1568.8

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H113").Value = 2513.8667
$ws.Range("I113").Value = 1770.8
$ws.Range("K113").Value = 1770.8
$ws.Range("M113").Value = 399.2

$ws.Range("H122").Value = 3141.5625
$ws.Range("I122").Value = 2400
$ws.Range("J122").Value = 3478.6365
$ws.Range("K122").Value = 7200
$ws.Range("L122").Value = 10435.9095
$ws.Range("M122").Value = -4750
$ws.Range("N122").Value = -15335.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1101.55
$ws.Range("I107").Value = 728.08
$ws.Range("J107").Value = 1724
$ws.Range("K107").Value = 2184.24
$ws.Range("L107").Value = 5172
$ws.Range("M107").Value = -264.2400000000002
$ws.Range("N107").Value = -9012

$ws.Range("H109").Value = 26800
$ws.Range("J109").Value = 26800
$ws.Range("L109").Value = 26800
$ws.Range("N109").Value = -29574

$ws.Range("H113").Value = 169172.83
$ws.Range("I113").Value = 333745.66
$ws.Range("J113").Value = 4600
$ws.Range("K113").Value = 1001236.98
$ws.Range("L113").Value = 13800
$ws.Range("M113").Value = -999066.98
$ws.Range("N113").Value = -18140

$ws.Range("H122").Value = 2494.963
$ws.Range("I122").Value = 2089.238
$ws.Range("J122").Value = 3915
$ws.Range("K122").Value = 6267.714
$ws.Range("L122").Value = 11745
$ws.Range("M122").Value = -3817.714
$ws.Range("N122").Value = -16645
